$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '26.203.43'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.80%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.678.28'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -1.03%  '

$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.76%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '211.93'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -3.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.5260'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -4.14%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.71%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2660'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -3.37%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06311'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -2.32%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '21.42'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.78%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07544'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.83%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.705.99'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +0.16%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.459'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.89%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.5657'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.23%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.000008042'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -3.68%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '66.62'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.49%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '26.259.82'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.70%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '1.003'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.77%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.834'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.02%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '188.41'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.78%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '10.43'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -5.08%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.210'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.76%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.74%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '149.19'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +0.25%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1253'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -5.75%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '7.600'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -4.07%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '16.02'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +1.36%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.06211'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -1.37%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.360'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.65%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.283'

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.500'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -3.08%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.443'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -4.47%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.635'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -2.99%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.004'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -3.88%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.6077'

$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.38%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.724'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +0.14%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.124'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.93%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.01616'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -2.47%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.084.83'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -2.89%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.8693'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -1.47%  '

$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -0.96%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '100.06'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.53%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.827.12'

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.00000000109'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.21%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '56.40'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -2.09%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -0.61%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.019'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -2.31%  '

$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.64%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.4255'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.13%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '5.987'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -1.99%  '
